# Update the week's dates on Sheet1 (row 5, columns B:H) to the next
# timesheet period (Feb 15 - Feb 21, 2021), and move the active selection
# from C11 to D6, matching the author's edit for "Added My Timesheets and
# New Folders".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = 44242
$ws.Range("C5").Value = 44243
$ws.Range("D5").Value = 44244
$ws.Range("E5").Value = 44245
$ws.Range("F5").Value = 44246
$ws.Range("G5").Value = 44247
$ws.Range("H5").Value = 44248

$ws.Range("D6").Select() | Out-Null
